$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 245
$ws1.Range("F3").Value = 247

# Update the same values on the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 245
$ws4.Range("F3").Value = 247
